$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.596.85'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.65%  '

# Row 3
$ws.Range("E3").Value = '  +2.04%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.21%  '

# Row 6
$ws.Range("E6").Value = '  +0.16%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3943'
$ws.Range("D7").ClearFormats()

# Row 8
$ws.Range("E8").Value = '  +2.40%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.527'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +6.37%  '

# Row 10
$ws.Range("E10").Value = '  +0.16%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.83'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +8.89%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08778'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.56%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.225'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +8.26%  '

# Row 14
$ws.Range("E14").Value = '  +2.78%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001325'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.84%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.618'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +5.45%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.696.65'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.10%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.10'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.68%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07056'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.69'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.49%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.870'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.11%  '

# Row 23
$ws.Range("E23").Value = '  +1.73%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.582.36'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.62%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.019'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.309'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.66%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.41'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.20%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.03'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.11%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.228'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.10%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.57'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.99%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.476'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +15.70%  '

# Row 32
$ws.Range("B32").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C32").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.882.10'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.07%  '

# Row 33
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.108'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.53%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.350'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +12.23%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08527'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.41'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.76%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.959'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.56%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2739'
$ws.Range("D38").ClearFormats()

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.56'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.77%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02775'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +10.44%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09081'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.48%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.464'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.44%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7707'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.47%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7193'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.74%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.46'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +4.46%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.540'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +6.31%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.218'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.92%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.354'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +12.64%  '

# Row 49
$ws.Range("E49").Value = '  +0.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.42'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.44%  '

# Row 51
$ws.Range("E51").Value = '  +3.51%  '
